$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Style cleanup: F6 was using a redundant duplicate style (fill explicitly
# "applied" but set to none, identical in appearance to the plain bordered/
# wrap-text style used elsewhere). Re-touching WrapText makes the engine
# resolve F6 back onto the shared/common style instead of the stale one.
$ws.Range("F6").WrapText = $true

# --- Re-label a few tiles (text only, same cell style) ---
# NOTE: order matters for shared-string append order, so Dream Temple is
# written before Hallowed Ground.
$ws.Range("D5").Value = "Dream Temple"
$ws.Range("D4").Value = "Hallowed Ground"
$ws.Range("D6").Value = "Hidden Path"

# --- Mark tiles as "locked" with an orange fill (new fill color) ---
$lockedColor = 49407   # RGB(255,192,0) -> BGR-packed OLE color
$ws.Range("D6").Interior.Color = $lockedColor
$ws.Range("G7").Interior.Color = $lockedColor
$ws.Range("H9").Interior.Color = $lockedColor
$ws.Range("D12").Interior.Color = $lockedColor
$ws.Range("K12").Interior.Color = $lockedColor

# --- New legend entries in column M (rows 3 and 4) ---
$startColor = 5296274  # RGB(146,208,80) -> BGR-packed OLE color (matches existing green)

$ws.Range("M4").Value = "Starting Tile"
$ws.Range("M4").Interior.Color = $startColor
$ws.Range("M4").Font.Size = 20

$ws.Range("M3").Value = "Locked - Quest needed to unlock area"
$ws.Range("M3").Interior.Color = $lockedColor
$ws.Range("M3").Font.Size = 20

# --- Move the active cell selection ---
$ws.Range("H9").Select() | Out-Null
